$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$excel.DisplayAlerts = $false

# ---------------------------------------------------------------------------
# 1. Un-merge the old "biorefinery" group (A9:A11). Row 9 stops being part of
#    that group (it becomes its own "Stream-natural gas" row), and the
#    "biorefinery" group shifts down to A10:A12 to make room for the new
#    "Price [USD/cf]" (natural-gas price) metric row.
# ---------------------------------------------------------------------------
$ws.Range("A9:A11").UnMerge()

$ws.Range("A9").Value  = "Stream-natural gas"
$ws.Range("A10").Value = "biorefinery"

# Give the new row 12 ("A12"/"B12") the same formatting as the rest of the
# table (copy down from row 11, which has the same column-A/column-B style).
$ws.Range("A11:F11").Copy()
$ws.Range("A12:F12").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Register the new merge A10:A12.
$ws.Range("A10:A12").MergeCells = $true
# Merging re-derives borders for the merged block and can introduce new
# (unused-elsewhere) border/style combinations; reapply the plain thin-box
# format (copied from A4, which carries the original, unmodified style) so
# the merged range keeps using the same style as before.
$ws.Range("A4").Copy()
$ws.Range("A10:A12").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 2. Update the "Parameter" labels in column B. A new metric
#    ("Price [USD/cf]") is inserted at row 9, pushing "Electricity price",
#    "Operating days" and "IRR" down by one row each; row 12 ends up with
#    "IRR [%]".
# ---------------------------------------------------------------------------
$ws.Range("B9").Value  = "Price [USD/cf]"
$ws.Range("B10").Value = "Electricity price [USD/kWh]"
$ws.Range("B11").Value = "Operating days [day/yr]"
$ws.Range("B12").Value = "IRR [%]"

# ---------------------------------------------------------------------------
# 3. Overwrite the numeric Spearman-correlation values for rows 4-12 with
#    their updated figures.
# ---------------------------------------------------------------------------
$ws.Range("C4").Value = 0.729936813936814
$ws.Range("D4").Value = 0.7945965025965027
$ws.Range("E4").Value = 0.776927432927433
$ws.Range("F4").Value = -0.9939877059877059

$ws.Range("C5").Value = 0.1514114594114594
$ws.Range("D5").Value = 0.06120524520524521
$ws.Range("E5").Value = 0.06061266061266062
$ws.Range("F5").Value = -0.02250547050547051

$ws.Range("C6").Value = 0.4820729060729061
$ws.Range("D6").Value = 0.5471493791493791
$ws.Range("E6").Value = 0.5697923577923578
$ws.Range("F6").Value = 0.1211601491601492

$ws.Range("C7").Value = 0.08503108903108904
$ws.Range("D7").Value = 0.04432156432156432
$ws.Range("E7").Value = 0.04442834042834043
$ws.Range("F7").Value = -0.03448671448671448

$ws.Range("C8").Value = -0.05147163947163947
$ws.Range("D8").Value = -0.04215713415713416
$ws.Range("E8").Value = -0.04291465891465892
$ws.Range("F8").Value = 0.02514200514200514

$ws.Range("C9").Value = 0.3675268275268275
$ws.Range("D9").Value = -0.01864434664434665
$ws.Range("E9").Value = -0.01982932382932383
$ws.Range("F9").Value = -0.0008835728835728836

$ws.Range("C10").Value = -0.005847797847797848
$ws.Range("D10").Value = -0.006176706176706178
$ws.Range("E10").Value = -0.006791718791718792
$ws.Range("F10").Value = 0.01357959757959758

$ws.Range("C11").Value = 0.05593863193863193
$ws.Range("D11").Value = 0.04639001839001839
$ws.Range("E11").Value = 0.04568619368619369
$ws.Range("F11").Value = -0.0670974790974791

$ws.Range("C12").Value = 0.03816187416187417
$ws.Range("D12").Value = 0.07230211230211231
$ws.Range("E12").Value = 0.07114871914871916
$ws.Range("F12").Value = -0.08064959664959664
